$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D ("Service"), shifting D..L to E..M
$ws.Range("D1").EntireColumn.Insert()

# Populate the new "Account number" column
$ws.Range("D1").Value = "Account number"
$ws.Range("D2").Value = "{d.meter[i].accountNumber}"
$ws.Range("D3").Value = "{d.meter[i + 1].accountNumber}"
